# The source data feed re-ordered two pairs of fixtures that share the same
# kickoff date/time (rows 25/26 and rows 88/89). The running index in column A
# stays tied to the row, but every other field (id, teams, scores, odds, ...)
# for the two fixtures is swapped.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowData {
    param($ws, [int]$row1, [int]$row2)

    $range1 = $ws.Range("B$row1:AD$row1")
    $range2 = $ws.Range("B$row2:AD$row2")

    $values1 = $range1.Value2
    $values2 = $range2.Value2

    $range1.Value2 = $values2
    $range2.Value2 = $values1
}

Swap-RowData $ws 25 26
Swap-RowData $ws 88 89
